$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D5").Value = 5378
$ws.Range("E5").Value = 8450
$ws.Range("F5").Value = 4391556
$ws.Range("G5").Value = "256-384"

$ws.Range("C14").Value = 1538
$ws.Range("D14").Value = 5378
$ws.Range("E14").Value = 8450
$ws.Range("F14").Value = 4391556
$ws.Range("G14").Value = "256-384"

$ws.Range("C17").Value = 71
$ws.Range("D17").Value = 75
$ws.Range("E17").Value = 79.400000000000006

$ws.Range("F16").Select()
